$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.448.15"
$ws.Range("E2").Value = "  +4.73%  "
$ws.Range("D3").Value = "2.747.40"
$ws.Range("E3").Value = "  +4.50%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'115.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.47%  "
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("D7").Value = "'0.538"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.32%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +5.51%  "
$ws.Range("D10").Value = "'41.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("D11").Value = "'0.0854"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.43%  "
$ws.Range("D12").Value = "'20.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("E14").Value = "  +4.82%  "
$ws.Range("D15").Value = "3.174.95"
$ws.Range("E15").Value = "  +4.34%  "
$ws.Range("D16").Value = "2.733.20"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("D17").Value = "'0.883"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").Value = "51.442.27"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("D19").Value = "'3.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.38%  "
$ws.Range("D20").Value = "'13.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.96%  "
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").Value = "0.0₃0972"
$ws.Range("E22").Value = "  +2.79%  "
$ws.Range("D23").Value = "'278.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.42%  "
$ws.Range("D24").Value = "'69.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("D26").Value = "'26.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").Value = "'2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").Value = "'35.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'49.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "'0.0825"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("D35").Value = "'19.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "'4.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").Value = "'127.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "'22.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("E42").Value = "  +7.71%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.113"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.70%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0343"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.03%  "
$ws.Range("D45").Value = "'2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.13%  "
$ws.Range("D46").Value = "2.087.31"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "'3.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("D48").Value = "'2.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.53%  "
$ws.Range("E49").Value = "  +6.17%  "
$ws.Range("D50").Value = "'8.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "'59.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.00%  "
